# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps for the
# 7de87ca3-4966-49b8-8cf0-c1197df20597 row (row 17) on both the
# zh-cn and de-de language sheets, reflecting a newer report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D17").Value = "2016-03-10 05:56:38"
$zhcn.Range("G17").Value = "2016-03-10 05:57:21"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D17").Value = "2016-03-10 05:56:47"
$dede.Range("G17").Value = "2016-03-10 05:57:37"
